$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Text = $old
    $find.Replacement.ClearFormatting()
    $find.Replacement.Text = $new
    $ok = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output ("WARNING: replace failed for: " + $old)
    }
    return $ok
}

# 1. Expand "Key/Value values can be arrays as well as key/value maps"
Replace-Text "Key/Value values can be arrays as well as key/value maps" "Key/Value pair values can be arrays or can be other key/value pairs" | Out-Null

# 2. Expand "Arrays can hold Key/Value Maps"
Replace-Text "Arrays can hold Key/Value Maps" "Arrays can hold Key/Value pairs" | Out-Null

# 3. Remove the stray mid-document "_GoBack" bookmark (this cascades the _Toc bookmark
#    ids down by one automatically, matching the renumbering seen in the diff).
try {
    $old_goback = $d.Bookmarks("_GoBack")
    $old_goback.Delete()
} catch {
    Write-Output ("No existing _GoBack bookmark to delete: " + $_.Exception.Message)
}

# 4. Update the platform hint text
Replace-Text "Hint: you will have to remove the VALID_PLATFORMS line from the make file (or add CYW943907*) " "Hint: you will have to remove the VALID_PLATFORMS line from the make file (or add WW101_*) " | Out-Null

# 5. Fix the capitalization of JSON_Parser -> JSON_parser in the exercise heading
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "(Advanced) Process a JSON document using " + [char]8220 + "JSON_Parser" + [char]8221
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "(Advanced) Process a JSON document using " + [char]8220 + "JSON_parser" + [char]8221
$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null

# 6. Update the final exercise description: "exercise (04)" -> "the previous exercise" and
#    JSON_Parser -> JSON_parser
Replace-Text "Write a program that will parse the same JSON document as exercise (04), but using the JSON_Parser library." "Write a program that will parse the same JSON document as the previous exercise, but using the JSON_parser library." | Out-Null

# 7. At the very end of the document there are two empty trailing paragraphs. Keep the
#    first empty, and put a fresh "_GoBack" bookmark (with the next free id) into the
#    second one, mirroring what Word does when it stamps the last-edit location back in.
$total = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($total)
$secondLastPara = $d.Paragraphs.Item($total - 1)
$insertPos = $lastPara.Range.Start
$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "done"
